# Generate Report for Handback
# The 473cc7fc-b280-4717-9643-499c8d7acb2a file has finished its handback
# round-trip for both locales, so flip its status from "Ready for handoff"
# to "Handed back: in sync with en-US" and stamp the handback datetimes.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B3").Value = "Handed back: in sync with en-US"
$overview.Range("C3").Value = "Handed back: in sync with en-US"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = "Handed back: in sync with en-US"
$zhcn.Range("H3").Value = "2016-03-21 00:40:53"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = "Handed back: in sync with en-US"
$dede.Range("H3").Value = "2016-03-21 00:40:59"
